$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E5").Value = 13.098
$ws.Range("E6").Value = 13.284
$ws.Range("D7").Value = -7.27
$ws.Range("B10").Value = 7.095000000000001
$ws.Range("B12").Value = 6.548999999999999
$ws.Range("C13").Value = -13.059
$ws.Range("B18").Value = 6.548999999999999
$ws.Range("D20").Value = -8.222
